# Apply the "verificar lo de fredy" edit to FMAQUINAS.xlsx
# - Update the "FECHA INICIAL" date text from 1/03/2021 to 1/04/2021
# - Fill in previously-zero "ACEITES Y FILTROS" (column D) expense values
#   for machines FOTON 1, FOTON 2, RETRO 5, RETRO 6 and SZN-114

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# FECHA INICIAL value (merged cell F6:I7).
# Assigning a date-like string straight to .Value makes Excel auto-convert it
# into a date serial number, so instead we push it in as a literal-text
# formula, then flatten the formula down to a static value via copy/paste
# (values only) and finally restore the merged-range formatting.
$ws.Range("F6").Formula = "=""1/04/2021"""
$ws.Range("F6").Copy()
$ws.Range("F6").PasteSpecial(-4163)
$ws.Range("F6:I7").Merge()

# ACEITES Y FILTROS column (D) values
$ws.Range("D9").Value = 45000    # FOTON 1
$ws.Range("D10").Value = 500000  # FOTON 2
$ws.Range("D15").Value = 50000   # RETRO 5
$ws.Range("D16").Value = 500000  # RETRO 6
$ws.Range("D17").Value = 500000  # SZN-114
